$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.48654696685571
$ws.Range("C2").Value = 8.240977258690494
$ws.Range("D2").Value = 9.235474447145386
$ws.Range("E2").Value = 13.51009371198972
$ws.Range("F2").Value = 31.73970379576936
$ws.Range("J2").Value = 9.944062962671381
$ws.Range("M2").Value = 16.71181436925981
$ws.Range("O2").Value = 23.80790824098548
$ws.Range("B3").Value = 13.91101951493271
$ws.Range("C3").Value = 7.730390033143972
$ws.Range("D3").Value = 9.216359259248561
$ws.Range("E3").Value = 13.53106918420202
$ws.Range("F3").Value = 31.83984451700819
$ws.Range("J3").Value = 9.974609070204281
$ws.Range("M3").Value = 16.52247938479741
$ws.Range("O3").Value = 23.91556842161778
$ws.Range("B4").Value = 13.54625247433378
$ws.Range("C4").Value = 7.397928264072198
$ws.Range("D4").Value = 9.205663590727967
$ws.Range("E4").Value = 13.54613840086089
$ws.Range("F4").Value = 31.91096345801167
$ws.Range("J4").Value = 9.994593257764556
$ws.Range("M4").Value = 16.40714083840774
$ws.Range("O4").Value = 23.98855442043285
$ws.Range("B5").Value = 13.39495098139633
$ws.Range("C5").Value = 7.257699424193824
$ws.Range("D5").Value = 9.201569485162482
$ws.Range("E5").Value = 13.55282988678279
$ws.Range("F5").Value = 31.94235665550615
$ws.Range("J5").Value = 10.00304635976545
$ws.Range("M5").Value = 16.36041412054752
$ws.Range("O5").Value = 24.02002025669163
$ws.Range("B6").Value = 13.36967364434069
$ws.Range("C6").Value = 7.234128484294719
$ws.Range("D6").Value = 9.200905709831332
$ws.Range("E6").Value = 13.55397425717083
$ws.Range("F6").Value = 31.94771481285352
$ws.Range("J6").Value = 10.00446868885881
$ws.Range("M6").Value = 16.35267307713041
$ws.Range("O6").Value = 24.02534901580603
$ws.Range("B7").Value = 13.54422242994908
$ws.Range("C7").Value = 7.396056282142316
$ws.Range("D7").Value = 9.205607301969509
$ws.Range("E7").Value = 13.54622641516094
$ws.Range("F7").Value = 31.91137708725528
$ws.Range("J7").Value = 9.994706006037065
$ws.Range("M7").Value = 16.40650949444533
$ws.Range("O7").Value = 23.98897181109146
$ws.Range("B8").Value = 14.29058675869822
$ws.Range("C8").Value = 8.068871208454977
$ws.Range("D8").Value = 9.228669086768724
$ws.Range("E8").Value = 13.51687145352206
$ws.Range("F8").Value = 31.77222708533204
$ws.Range("J8").Value = 9.954340474956448
$ws.Range("M8").Value = 16.64636901303059
$ws.Range("O8").Value = 23.84359681464534
$ws.Range("B9").Value = 15.65587662853109
$ws.Range("C9").Value = 9.237314822502167
$ws.Range("D9").Value = 9.282031425263806
$ws.Range("E9").Value = 13.47668615368209
$ws.Range("F9").Value = 31.57621933887837
$ws.Range("J9").Value = 9.884915177762675
$ws.Range("M9").Value = 17.12200196913671
$ws.Range("O9").Value = 23.61343776194123
$ws.Range("B10").Value = 16.58992355463615
$ws.Range("C10").Value = 10.00344169797895
$ws.Range("D10").Value = 9.326025543748607
$ws.Range("E10").Value = 13.45775577746566
$ws.Range("F10").Value = 31.47962799472345
$ws.Range("J10").Value = 9.839815672502482
$ws.Range("M10").Value = 17.47200939203139
$ws.Range("O10").Value = 23.47824134432585
$ws.Range("B11").Value = 16.99833506246387
$ws.Range("C11").Value = 10.33193769106958
$ws.Range("D11").Value = 9.347039575943537
$ws.Range("E11").Value = 13.45144205011695
$ws.Range("F11").Value = 31.44608545964149
$ws.Range("J11").Value = 9.820576209307823
$ws.Range("M11").Value = 17.63080181295798
$ws.Range("O11").Value = 23.42418512782371
$ws.Range("B12").Value = 17.15051473763639
$ws.Range("C12").Value = 10.45345570203343
$ws.Range("D12").Value = 9.355137320339441
$ws.Range("E12").Value = 13.44938124640451
$ws.Range("F12").Value = 31.43488605235892
$ws.Range("J12").Value = 9.813473917002234
$ws.Range("M12").Value = 17.69082512370645
$ws.Range("O12").Value = 23.40479295468574
$ws.Range("B13").Value = 17.11785180156797
$ws.Range("C13").Value = 10.42741248761481
$ws.Range("D13").Value = 9.353387152423815
$ws.Range("E13").Value = 13.44981040427002
$ws.Range("F13").Value = 31.43723111780236
$ws.Range("J13").Value = 9.814995376027476
$ws.Range("M13").Value = 17.67790355863963
$ws.Range("O13").Value = 23.40892136226829
$ws.Range("B14").Value = 17.01090516702175
$ws.Range("C14").Value = 10.34199267495072
$ws.Range("D14").Value = 9.347702994628056
$ws.Range("E14").Value = 13.45126589390841
$ws.Range("F14").Value = 31.44513392826691
$ws.Range("J14").Value = 9.819988227965259
$ws.Range("M14").Value = 17.63574239236164
$ws.Range("O14").Value = 23.42256808251957
$ws.Range("B15").Value = 16.94507182906415
$ws.Range("C15").Value = 10.28929610242234
$ws.Range("D15").Value = 9.344239426847194
$ws.Range("E15").Value = 13.45220039560199
$ws.Range("F15").Value = 31.45017048893863
$ws.Range("J15").Value = 9.82307035332823
$ws.Range("M15").Value = 17.60990204365361
$ws.Range("O15").Value = 23.43106765410072
$ws.Range("B16").Value = 16.56289168773323
$ws.Range("C16").Value = 9.981571036366315
$ws.Range("D16").Value = 9.324672051313268
$ws.Range("E16").Value = 13.45821460426529
$ws.Range("F16").Value = 31.48202994358
$ws.Range("J16").Value = 9.841098678834129
$ws.Range("M16").Value = 17.46161933630919
$ws.Range("O16").Value = 23.48192447715985
$ws.Range("B17").Value = 16.32413275469785
$ws.Range("C17").Value = 9.787665353699913
$ws.Range("D17").Value = 9.31292183631775
$ws.Range("E17").Value = 13.46249243060564
$ws.Range("F17").Value = 31.50424317059788
$ws.Range("J17").Value = 9.852485214181854
$ws.Range("M17").Value = 17.37050981741293
$ws.Range("O17").Value = 23.51503554513971
$ws.Range("B18").Value = 16.1852579700302
$ws.Range("C18").Value = 9.674251354537729
$ws.Range("D18").Value = 9.30625786478679
$ws.Range("E18").Value = 13.46516923046714
$ws.Range("F18").Value = 31.51799789787254
$ws.Range("J18").Value = 9.859154617685441
$ws.Range("M18").Value = 17.31806849659501
$ws.Range("O18").Value = 23.5347805310983
$ws.Range("B19").Value = 16.13797513107634
$ws.Range("C19").Value = 9.635527607142027
$ws.Range("D19").Value = 9.304017890508534
$ws.Range("E19").Value = 13.46611270886407
$ws.Range("F19").Value = 31.52282279364503
$ws.Range("J19").Value = 9.861433411280968
$ws.Range("M19").Value = 17.30030775003258
$ws.Range("O19").Value = 23.5415859139313
$ws.Range("B20").Value = 16.34970996828608
$ws.Range("C20").Value = 9.808501980708433
$ws.Range("D20").Value = 9.314162919081662
$ws.Range("E20").Value = 13.46201466390617
$ws.Range("F20").Value = 31.50177723201544
$ws.Range("J20").Value = 9.851260663070658
$ws.Range("M20").Value = 17.38021280540426
$ws.Range("O20").Value = 23.51143827692524
$ws.Range("B21").Value = 17.04238597387436
$ws.Range("C21").Value = 10.36716057674707
$ws.Range("D21").Value = 9.349368795447159
$ws.Range("E21").Value = 13.45082942672123
$ws.Range("F21").Value = 31.44277185223563
$ws.Range("J21").Value = 9.818516734550823
$ws.Range("M21").Value = 17.64812944320007
$ws.Range("O21").Value = 23.41853040472588
$ws.Range("B22").Value = 17.48061394849053
$ws.Range("C22").Value = 10.71551353248104
$ws.Range("D22").Value = 9.373193151523459
$ws.Range("E22").Value = 13.44544287496419
$ws.Range("F22").Value = 31.41296844032321
$ws.Range("J22").Value = 9.798184822813537
$ws.Range("M22").Value = 17.82257854905934
$ws.Range("O22").Value = 23.36409479198404
$ws.Range("B23").Value = 17.24807897779555
$ws.Range("C23").Value = 10.53112347937823
$ws.Range("D23").Value = 9.360404317873446
$ws.Range("E23").Value = 13.44814190968994
$ws.Range("F23").Value = 31.42807136292397
$ws.Range("J23").Value = 9.808938705060109
$ws.Range("M23").Value = 17.72954631502131
$ws.Range("O23").Value = 23.39257072701384
$ws.Range("B24").Value = 16.33815151383497
$ws.Range("C24").Value = 9.799087769089427
$ws.Range("D24").Value = 9.313601540646205
$ws.Range("E24").Value = 13.46222998496326
$ws.Range("F24").Value = 31.50288901811567
$ws.Range("J24").Value = 9.851813899049638
$ws.Range("M24").Value = 17.37582627174573
$ws.Range("O24").Value = 23.51306239390518
$ws.Range("B25").Value = 15.29807495653765
$ws.Range("C25").Value = 8.937427633310714
$ws.Range("D25").Value = 9.266741412232815
$ws.Range("E25").Value = 13.48569647954035
$ws.Range("F25").Value = 31.6209552838721
$ws.Range("J25").Value = 9.902657390386672
$ws.Range("M25").Value = 16.99304129737314
$ws.Range("O25").Value = 23.66977821040415
